$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.124.29'
$ws.Cells.Item(2, 5).Value = '  +0.14%  '
$ws.Cells.Item(3, 4).Value = '1.665.64'
$ws.Cells.Item(3, 5).Value = '  -0.66%  '
$ws.Cells.Item(4, 5).Value = '  -0.15%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '209.62'
$ws.Cells.Item(5, 5).Value = '  -0.63%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.5207'
$ws.Cells.Item(6, 5).Value = '  -1.49%  '
$ws.Cells.Item(7, 5).Value = '  -0.12%  '
$ws.Cells.Item(8, 5).Value = '  -2.98%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.06318'
$ws.Cells.Item(9, 5).Value = '  +0.07%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '21.02'
$ws.Cells.Item(10, 5).Value = '  -1.24%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07538'
$ws.Cells.Item(11, 5).Value = '  -0.52%  '
$ws.Cells.Item(12, 4).Value = '1.670.75'
$ws.Cells.Item(12, 5).Value = '  -0.44%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '4.405'
$ws.Cells.Item(13, 5).Value = '  -2.27%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.5410'
$ws.Cells.Item(14, 5).Value = '  -4.82%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.000007986'
$ws.Cells.Item(15, 5).Value = '  -1.92%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '66.26'
$ws.Cells.Item(16, 5).Value = '  +0.88%  '
$ws.Cells.Item(17, 4).Value = '26.161.20'
$ws.Cells.Item(17, 5).Value = '  +0.11%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '1.003'
$ws.Cells.Item(18, 5).Value = '  -0.18%  '
$ws.Cells.Item(19, 5).Value = '  -2.78%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '186.86'
$ws.Cells.Item(20, 5).Value = '  -1.29%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '10.23'
$ws.Cells.Item(21, 5).Value = '  -3.64%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.214'
$ws.Cells.Item(22, 5).Value = '  +0.12%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.004'
$ws.Cells.Item(23, 5).Value = '  -0.08%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '149.65'
$ws.Cells.Item(24, 5).Value = '  +0.72%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.1236'
$ws.Cells.Item(25, 5).Value = '  -1.60%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '7.397'
$ws.Cells.Item(26, 5).Value = '  -3.34%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '15.70'
$ws.Cells.Item(27, 5).Value = '  -2.19%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.06266'
$ws.Cells.Item(28, 5).Value = '  -1.77%  '
$ws.Cells.Item(29, 5).Value = '  +0.62%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.274'
$ws.Cells.Item(30, 5).Value = '  -0.95%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '3.489'
$ws.Cells.Item(31, 5).Value = '  -1.44%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.400'
$ws.Cells.Item(32, 5).Value = '  -3.93%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.635'
$ws.Cells.Item(33, 5).Value = '  -2.33%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.9977'
$ws.Cells.Item(34, 5).Value = '  -1.29%  '
$ws.Cells.Item(35, 2).Value = 'HuobiToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '2.395'
$ws.Cells.Item(35, 5).Value = '  -0.90%  '
$ws.Cells.Item(36, 2).Value = 'MXToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.758'
$ws.Cells.Item(36, 5).Value = '  +1.32%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.5961'
$ws.Cells.Item(37, 5).Value = '  -1.74%  '
$ws.Cells.Item(38, 4).Value = '1.108.87'
$ws.Cells.Item(38, 5).Value = '  +1.05%  '
$ws.Cells.Item(39, 5).Value = '  -0.43%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '6.048'
$ws.Cells.Item(40, 5).Value = '  -1.71%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.8624'
$ws.Cells.Item(41, 5).Value = '  -0.97%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.004'
$ws.Cells.Item(42, 5).Value = '  -0.09%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '100.59'
$ws.Cells.Item(43, 5).Value = '  +0.52%  '
$ws.Cells.Item(44, 4).Value = '1.813.88'
$ws.Cells.Item(44, 5).Value = '  -0.78%  '
$ws.Cells.Item(45, 5).Value = '  -0.99%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '55.23'
$ws.Cells.Item(46, 5).Value = '  -3.25%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.001'
$ws.Cells.Item(47, 5).Value = '  -0.61%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '8.059'
$ws.Cells.Item(48, 5).Value = '  +0.59%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.05243'
$ws.Cells.Item(49, 5).Value = '  -0.20%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.4235'
$ws.Cells.Item(50, 5).Value = '  -0.75%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '5.867'
$ws.Cells.Item(51, 5).Value = '  -1.53%  '
